$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (the surviving "InsuranceQuoteByDB Select" row) ---

# E2: StepInfo text changes, and gains wrap text
$ws.Cells.Item(2,5).Value = "get insurance quote id;for insurance premium;query;information on the;insurance premium;"
$ws.Cells.Item(2,5).WrapText = $true

# N2: Csvson column gains the CSV payload (previously blank), with wrap text
$n2Text = "id,insurance_premium_amount, insurance_premium_currency, policy_limit_amount" + [char]10 + "i~[request_Id],d~500.00,CHF,d~50000.00"
$ws.Cells.Item(2,14).Value = $n2Text
$ws.Cells.Item(2,14).WrapText = $true

# --- Remove the old "InsuranceQuoteByDB Verify" row (row 3) entirely ---
$ws.Rows.Item(3).Delete()

# --- Column E width update ---
$ws.Columns.Item(5).ColumnWidth = 42.15

# --- Update the active selection / view to match the final state ---
$ws.Activate() | Out-Null
$ws.Range("E3").Select() | Out-Null
